{"js": "// Replace every occurrence of \"ISNSE\" with \"Internet as a Single Computer\"\n// (the acronym is spelled out), keeping the existing run formatting\n// (orange/FFC000 color) intact, and relocate the \"_GoBack\" last-edit\n// bookmark the way Word itself would after these edits.\n\nconst body = context.document.body;\n\nconst results = body.search(\"ISNSE\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Internet as a Single Computer\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark (Word's \"last edit location\" marker) so that it\n// sits right after the paragraph ending in \"...Most user interfaces are\",\n// matching where Word leaves it once the edits above have been made.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst anchor = body.search(\"Most user interfaces are\", { matchCase: true });\nanchor.load(\"text\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const end = anchor.items[0].getRange(Word.RangeLocation.end);\n  end.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Spell out the \"ISNSE\" acronym as \"Internet as a Single Computer\" everywhere\n# it appears in the document, keeping the existing (orange / FFC000) run\n# formatting intact, then relocate the \"_GoBack\" last-edit bookmark the way\n# Word itself leaves it after these edits (right after the paragraph that\n# ends in \"...Most user interfaces are\").\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"ISNSE\"\n$find.Replacement.Text = \"Internet as a Single Computer\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Remove the existing \"_GoBack\" bookmark, wherever Word currently has it.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-anchor \"_GoBack\" right after \"...Most user interfaces are\" \u2014 delete and\n# retype the final character of that run so the newly added bookmark\n# collapses to a zero-length marker at that exact point (matching how Word\n# leaves _GoBack after an edit).\n$anchorFind = $d.Content.Find\n$anchorFind.ClearFormatting()\n$anchorFind.Text = \"Most user interfaces are\"\n$anchorFound = $anchorFind.Execute()\n\nif ($anchorFound) {\n  $anchorRange = $anchorFind.Parent\n  $lastChar = $d.Range($anchorRange.End - 1, $anchorRange.End)\n  $lastCharText = $lastChar.Text\n  $d.Bookmarks.Add(\"_GoBack\", $lastChar) | Out-Null\n  $lastChar.Delete()\n  $d.Range($anchorRange.End - 1, $anchorRange.End - 1).InsertAfter($lastCharText) | Out-Null\n}\n"}
